$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 already carries the bold/bordered/centered label style used for all
# column-A entries; clone that formatting onto the newly added A3:A5 cells
# (copy + paste-special-formats keeps the existing style index instead of
# minting a near-duplicate one).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

# Now write the new values for all five data rows.
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 94

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 93

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 92

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 36
